$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the "Objectives:" row (row 11), pushing
# "Programa resumido:" and everything below it down by two rows.
$ws.Rows("12:13").Insert()

# New label row (column A only), mirrors the style of the other label-only
# rows (e.g. old row 16 "Avaliação:").
$ws.Range("A12").Value = "Docentes responsáveis:"

# New value row (columns B and C only, same text in both, mirrors e.g. the
# "LOB1012" row at the bottom of the sheet).
$ws.Range("B13").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("C13").Value = "8855158 - Morun Bernardino Neto"

# The row Insert() stamped every column of the new rows with the column's
# default format; drop the format+content on the cells that should stay
# empty so the two rows only carry cells in the columns actually used
# (matches the other label-only / value-only rows elsewhere in the sheet).
$ws.Range("B12:C12").Clear()
$ws.Range("A13").Clear()
